$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rows 99-106 in column D: these FreeCodeCamp lessons are now complete.
# Row 99 text was already updated to "Passed..." previously but still carried
# the "Not Passed" (orange) highlight style - just clear that leftover fill.
# Rows 100-106 get both their text flipped from "Not Passed..." to "Passed..."
# and the leftover highlight cleared.

$ws.Range("D100").Value = "PassedNesting For Loops"
$ws.Range("D101").Value = "PassedIterate with JavaScript Do...While Loops"
$ws.Range("D102").Value = "PassedProfile Lookup"
$ws.Range("D103").Value = "PassedGenerate Random Fractions with JavaScript"
$ws.Range("D104").Value = "PassedGenerate Random Whole Numbers with JavaScript"
$ws.Range("D105").Value = "PassedGenerate Random Whole Numbers within a Range"
$ws.Range("D106").Value = "PassedUse the parseInt Function"

# Clear the stale "Not Passed" highlight fill on D99:D106 now that all of
# these rows reflect completed ("Passed") lessons.
$ws.Range("D99:D106").ClearFormats()

# Move the selection / scroll position to the rows just updated.
$ws.Range("D98:D106").Select()
